$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11224.223
$ws.Range("I6").Value = 12584.25
$ws.Range("K6").Value = 37752.75
$ws.Range("M6").Value = -37640.75
$ws.Range("H70").Value = 50002772
$ws.Range("I70").Value = 2333.3333
$ws.Range("J70").Value = 58826380
$ws.Range("K70").Value = 6999.999899999999
$ws.Range("L70").Value = 176479140
$ws.Range("M70").Value = -6729.999899999999
$ws.Range("N70").Value = -176479680
$ws.Range("H73").Value = 50002772
$ws.Range("I73").Value = 2333.3333
$ws.Range("J73").Value = 58826380
$ws.Range("K73").Value = 6999.999899999999
$ws.Range("L73").Value = 176479140
$ws.Range("M73").Value = -6063.999899999999
$ws.Range("N73").Value = -176481012
$ws.Range("H87").Value = 333366660
$ws.Range("J87").Value = 333366660
$ws.Range("L87").Value = 333366660
$ws.Range("N87").Value = -333369156
$ws.Range("H90").Value = 333366660
$ws.Range("J90").Value = 333366660
$ws.Range("L90").Value = 1000099980
$ws.Range("N90").Value = -1000112460
$ws.Range("H92").Value = 2334.9092
$ws.Range("I92").Value = 3995
$ws.Range("J92").Value = 1386.2858
$ws.Range("K92").Value = 3995
$ws.Range("L92").Value = 1386.2858
$ws.Range("M92").Value = -2747
$ws.Range("N92").Value = -3882.2858
$ws.Range("H113").Value = 3682.8333
$ws.Range("I113").Value = 3595.3333
$ws.Range("J113").Value = 4032.8333
$ws.Range("K113").Value = 3595.3333
$ws.Range("L113").Value = 4032.8333
$ws.Range("M113").Value = -341.3332999999998
$ws.Range("N113").Value = -10540.8333
$ws.Range("H125").Value = 30469.637
$ws.Range("J125").Value = 28966.75
$ws.Range("L125").Value = 260700.75
$ws.Range("N125").Value = -265620.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1343.625
$ws.Range("J4").Value = 2079.8
$ws.Range("L4").Value = 2079.8
$ws.Range("N4").Value = -2311.8
$ws.Range("H6").Value = 404166.66
$ws.Range("J6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("N6").Value = -15346
$ws.Range("H63").Value = 12572
$ws.Range("I63").Value = 12668
$ws.Range("K63").Value = 12668
$ws.Range("M63").Value = -11982
$ws.Range("H66").Value = 12572
$ws.Range("I66").Value = 12668
$ws.Range("K66").Value = 63340
$ws.Range("M66").Value = -59908
$ws.Range("H74").Value = 2679.319
$ws.Range("I74").Value = 2280.9666
$ws.Range("J74").Value = 3382.2942
$ws.Range("K74").Value = 2280.9666
$ws.Range("L74").Value = 3382.2942
$ws.Range("M74").Value = -1406.9666
$ws.Range("N74").Value = -5130.2942
$ws.Range("H77").Value = 2679.319
$ws.Range("I77").Value = 2280.9666
$ws.Range("J77").Value = 3382.2942
$ws.Range("K77").Value = 11404.833
$ws.Range("L77").Value = 16911.471
$ws.Range("M77").Value = -7036.833000000001
$ws.Range("N77").Value = -25647.471
$ws.Range("H122").Value = 3004833.5
$ws.Range("I122").Value = 4630859.5
$ws.Range("J122").Value = 2939.2307
$ws.Range("K122").Value = 13892578.5
$ws.Range("L122").Value = 8817.6921
$ws.Range("M122").Value = -13890128.5
$ws.Range("N122").Value = -13717.6921
$ws.Range("H132").Value = 1505.683
$ws.Range("I132").Value = 1505.683
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4517.049
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -1987.049

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 781.875
$ws.Range("I80").Value = 569.8333
$ws.Range("K80").Value = 569.8333
$ws.Range("M80").Value = 428.1667
$ws.Range("H83").Value = 781.875
$ws.Range("I83").Value = 569.8333
$ws.Range("K83").Value = 2849.1665
$ws.Range("M83").Value = 2142.8335
$ws.Range("H107").Value = 3218.1333
$ws.Range("I107").Value = 2013.9459
$ws.Range("J107").Value = 8787.5
$ws.Range("K107").Value = 2013.9459
$ws.Range("L107").Value = 8787.5
$ws.Range("M107").Value = -93.94589999999994
$ws.Range("N107").Value = -12627.5
$ws.Range("H115").Value = 49750
$ws.Range("J115").Value = 49750
$ws.Range("L115").Value = 49750
$ws.Range("N115").Value = -52884
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1173.1052
$ws.Range("I4").Value = 1788.3334
$ws.Range("J4").Value = 619.4
$ws.Range("K4").Value = 1788.3334
$ws.Range("L4").Value = 619.4
$ws.Range("M4").Value = -1676.3334
$ws.Range("N4").Value = -843.4
$ws.Range("H5").Value = 229.125
$ws.Range("J5").Value = 499
$ws.Range("L5").Value = 499
$ws.Range("N5").Value = -723
$ws.Range("H16").Value = 1640.4166
$ws.Range("I16").Value = 1318.6
$ws.Range("J16").Value = 3249.5
$ws.Range("K16").Value = 1318.6
$ws.Range("L16").Value = 3249.5
$ws.Range("M16").Value = -1031.6
$ws.Range("N16").Value = -3823.5
$ws.Range("H92").Value = 62316.332
$ws.Range("J92").Value = 62316.332
$ws.Range("L92").Value = 62316.332
$ws.Range("N92").Value = -67308.33199999999
$ws.Range("H107").Value = 738.5714
$ws.Range("I107").Value = 418.46667
$ws.Range("J107").Value = 978.65
$ws.Range("K107").Value = 418.46667
$ws.Range("L107").Value = 978.65
$ws.Range("M107").Value = 1501.53333
$ws.Range("N107").Value = -4818.65
$ws.Range("H113").Value = 1640.4166
$ws.Range("I113").Value = 1318.6
$ws.Range("J113").Value = 3249.5
$ws.Range("K113").Value = 1318.6
$ws.Range("L113").Value = 3249.5
$ws.Range("M113").Value = 851.4000000000001
$ws.Range("N113").Value = -7589.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25689538
$ws.Range("I4").Value = 27390320
$ws.Range("J4").Value = 17428600
$ws.Range("K4").Value = 82170960
$ws.Range("L4").Value = 52285800
$ws.Range("M4").Value = -82170848
$ws.Range("N4").Value = -52286024
$ws.Range("H94").Value = 7708.1665
$ws.Range("I94").Value = 2666.6667
$ws.Range("K94").Value = 8000.000100000001
$ws.Range("M94").Value = -7324.000100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1167.7142
$ws.Range("I3").Value = 347.5
$ws.Range("J3").Value = 1495.8
$ws.Range("K3").Value = 347.5
$ws.Range("L3").Value = 1495.8
$ws.Range("M3").Value = -231.5
$ws.Range("N3").Value = -1727.8
$ws.Range("H11").Value = 2254357.2
$ws.Range("I11").Value = 3835876.8
$ws.Range("J11").Value = 57802.39
$ws.Range("K11").Value = 3835876.8
$ws.Range("L11").Value = 57802.39
$ws.Range("M11").Value = -3835737.8
$ws.Range("N11").Value = -58080.39
$ws.Range("H132").Value = 62501704
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 71433360
$ws.Range("I82").Value = 1110
$ws.Range("J82").Value = 111117950
$ws.Range("K82").Value = 1110
$ws.Range("L82").Value = 111117950
$ws.Range("M82").Value = -749
$ws.Range("N82").Value = -111118672
$ws.Range("H85").Value = 71433360
$ws.Range("I85").Value = 1110
$ws.Range("J85").Value = 111117950
$ws.Range("K85").Value = 1110
$ws.Range("L85").Value = 111117950
$ws.Range("M85").Value = 138
$ws.Range("N85").Value = -111120446

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("N11").Value = 0
$ws.Range("H107").Value = 29412804
$ws.Range("I107").Value = 752.6667
$ws.Range("K107").Value = 2258.0001
$ws.Range("M107").Value = -338.0001000000002
